$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Roll back the student insert buffer values
$ws.Range("A2").Value = 1150
$ws.Range("N2").Value = 1

# Restore selection to C7
$ws.Range("C7").Select()
